$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last two data rows (old rows 6 and 7) are removed entirely; the
# remaining rows (2-5) keep their row numbers and the used range shrinks
# to A1:T5.
$ws.Rows("6:7").Delete() | Out-Null

# Row 2 (Inflammatory-Mac / Npy / Npy1r / ECs): columns G:T get new TPM-based values.
$row2 = @(0.4700283333333333, 1.410085, 0.1029504401632623, 0.1029504401632623, 1, 0.3333333333333333, 0.1289933333333333, 0.38698, 0.05394679539474087, 0.0787948288373843, 0.06063052147777777, 0.5456746933, 0.005553846331286022, 0.008111962311397623)
$col = 7
foreach ($v in $row2) {
  $ws.Cells.Item(2, $col).Value = $v
  $col = $col + 1
}

# Row 3 (Inflammatory-Mac / Npy / Npy1r / MuSCs): columns G:T get new TPM-based values.
$row3 = @(0.4700283333333333, 1.410085, 0.1029504401632623, 0.1029504401632623, 2, 1, 2.262128, 4.524256, 0.9460532046052591, 0.9212051711626156, 1.063264253626667, 6.37958552176, 0.09739659383197624, 0.09483847785186462)
$col = 7
foreach ($v in $row3) {
  $ws.Cells.Item(3, $col).Value = $v
  $col = $col + 1
}

# Row 4: sending cluster changes from Inflammatory-Mac to Resolving-Mac, and the
# target cluster changes from Resolving-Mac to ECs; columns E:T get new values.
$ws.Cells.Item(4, 1).Value = "Resolving-Mac"
$ws.Cells.Item(4, 4).Value = "ECs"
$row4 = @(3, 1, 4.095550333333333, 12.286651, 0.8970495598367377, 0.8970495598367377, 1, 0.3333333333333333, 0.1289933333333333, 0.38698, 0.05394679539474087, 0.0787948288373843, 0.528298689331111, 4.75468820398, 0.04839294906345485, 0.07068286652598667)
$col = 5
foreach ($v in $row4) {
  $ws.Cells.Item(4, $col).Value = $v
  $col = $col + 1
}

# Row 5: target cluster changes from ECs to MuSCs; columns E:T get new values.
$ws.Cells.Item(5, 4).Value = "MuSCs"
$row5 = @(3, 1, 4.095550333333333, 12.286651, 0.8970495598367377, 0.8970495598367377, 2, 1, 2.262128, 4.524256, 0.9460532046052591, 0.9212051711626156, 9.264659084442666, 55.587954506656, 0.8486566107732828, 0.826366693310751)
$col = 5
foreach ($v in $row5) {
  $ws.Cells.Item(5, $col).Value = $v
  $col = $col + 1
}
